# Updated cryptos list on Sat Jun 29 13:25:28 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto
# ranking table on Sheet1. Cells D/E are stored as plain text in the
# workbook (percentages are padded with spaces, and many "prices" use a
# dotted thousands-separator format like "61.047.02" that must stay
# textual). For price cells whose new value would otherwise be
# auto-parsed by Excel as a genuine number (dropping trailing zeros,
# e.g. "1.00" -> 1), the cell is pre-formatted as Text ("@") so the
# literal string is preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.047.02"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "3.392.90"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.43"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.06"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D7").Value = "3.392.09"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.63"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("E11").Value = "  -2.70%  "

$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("D13").Value = "3.971.59"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.126"
$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.88"
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "3.386.67"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").Value = "61.086.13"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("E20").Value = "  -4.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.90"
$ws.Range("E21").Value = "  -4.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.55"
$ws.Range("E22").Value = "  -4.30%  "

$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.59"
$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -4.29%  "

$ws.Range("D27").Value = "3.531.21"
$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.34"
$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("E32").Value = "  -1.10%  "

$ws.Range("E33").Value = "  -6.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.47"
$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.86"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "3.425.33"
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.01"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("E40").Value = "  -4.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0770"
$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.21"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("E45").Value = "  -2.48%  "

$ws.Range("E46").Value = "  -3.26%  "

$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("D48").Value = "2.475.74"
$ws.Range("E48").Value = "  -5.48%  "

$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.90"
$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0266"
$ws.Range("E51").Value = "  +1.42%  "
